$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.161.32"
$ws.Range("E2").Value = "  +5.89%  "

$ws.Range("D3").Value = "1.917.01"
$ws.Range("E3").Value = "  +2.59%  "

$ws.Range("E4").Value = "  -0.91%  "

$ws.Range("D5").Value = "329.72"
$ws.Range("E5").Value = "  +4.55%  "

$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("D7").Value = "0.5237"
$ws.Range("E7").Value = "  +3.17%  "

$ws.Range("D8").Value = "0.4081"
$ws.Range("E8").Value = "  +4.43%  "

$ws.Range("D9").Value = "0.08523"
$ws.Range("E9").Value = "  +2.05%  "

$ws.Range("D10").Value = "42.86"
$ws.Range("E10").Value = "  +0.96%  "

$ws.Range("D11").Value = "1.123"
$ws.Range("E11").Value = "  +1.71%  "

$ws.Range("D12").Value = "22.44"
$ws.Range("E12").Value = "  +10.40%  "

$ws.Range("D13").Value = "6.453"
$ws.Range("E13").Value = "  +4.08%  "

$ws.Range("D14").Value = "1.924.72"
$ws.Range("E14").Value = "  +3.24%  "

$ws.Range("D15").Value = "7.389"
$ws.Range("E15").Value = "  +2.08%  "

$ws.Range("E16").Value = "  -0.95%  "

$ws.Range("D17").Value = "95.00"
$ws.Range("E17").Value = "  +4.08%  "

$ws.Range("D18").Value = "0.00001114"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").Value = "0.06695"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").Value = "18.38"
$ws.Range("E20").Value = "  +4.21%  "

$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("D22").Value = "6.008"
$ws.Range("E22").Value = "  +1.84%  "

$ws.Range("D23").Value = "30.177.14"
$ws.Range("E23").Value = "  +5.72%  "

$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +2.19%  "

$ws.Range("D25").Value = "2.215"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("D26").Value = "2.136.20"
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("D27").Value = "160.70"
$ws.Range("E27").Value = "  +2.32%  "

$ws.Range("D28").Value = "21.13"
$ws.Range("E28").Value = "  +2.91%  "

$ws.Range("D29").Value = "2.414"
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("D30").Value = "128.91"
$ws.Range("E30").Value = "  +2.42%  "

$ws.Range("D31").Value = "1.081"
$ws.Range("E31").Value = "  +4.16%  "

$ws.Range("D32").Value = "0.1064"
$ws.Range("E32").Value = "  +2.41%  "

$ws.Range("D33").Value = "5.992"
$ws.Range("E33").Value = "  +4.22%  "

$ws.Range("D34").Value = "3.641"
$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("D35").Value = "0.02490"
$ws.Range("E35").Value = "  +1.40%  "

$ws.Range("D36").Value = "0.06591"
$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("D37").Value = "0.2205"
$ws.Range("E37").Value = "  +2.15%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.229"
$ws.Range("E38").Value = "  +4.20%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.175"
$ws.Range("E39").Value = "  +2.77%  "

$ws.Range("D40").Value = "8.885"
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").Value = "0.6532"
$ws.Range("E41").Value = "  +2.76%  "

$ws.Range("D42").Value = "11.63"
$ws.Range("E42").Value = "  +4.98%  "

$ws.Range("D43").Value = "1.241"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "0.6145"
$ws.Range("E44").Value = "  +2.50%  "

$ws.Range("D45").Value = "13.26"
$ws.Range("E45").Value = "  +1.37%  "

$ws.Range("D46").Value = "3.748"
$ws.Range("E46").Value = "  +1.82%  "

$ws.Range("D47").Value = "2.084"
$ws.Range("E47").Value = "  +4.33%  "

$ws.Range("E48").Value = "  +2.77%  "

$ws.Range("D49").Value = "124.58"
$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("D50").Value = "1.167"
$ws.Range("E50").Value = "  +3.69%  "

$ws.Range("D51").Value = "79.67"
$ws.Range("E51").Value = "  +4.48%  "
